$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Sheet name: $($ws.Name)"
$v = $ws.Range("D18").Value
Write-Host "D18 type: $($v.GetType().Name)"
Write-Host "D18 value: $v"
$v2 = $ws.Range("A2").Value
Write-Host "A2 value: $v2"
$v3 = $ws.Cells.Item(2, 1).Value
Write-Host "Cells(2,1) value: $v3"
